$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# Note: the ColumnWidth COM property is offset from the raw OOXML "width"
# attribute by 5/6 (0.8333...), so subtract that to land on the exact
# target width stored in the XML.
$ws.Columns.Item(1).ColumnWidth = 73 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 61 - (5/6)

# Update cell values
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/configureTestEnvironmentAndRun-test-data"
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTest-test-data"
